# Auto-generated Excel COM-interop script to apply numeric corrections
# to the Kujata_Profits workbook (currentAveragePrice / Leve price columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 850
$ws.Range("J13").Value = 1200
$ws.Range("L13").Value = 1200
$ws.Range("N13").Value = -1538
# Row 33
$ws.Range("H33").Value = 409.03125
$ws.Range("I33").Value = 379.84
$ws.Range("J33").Value = 513.2857
$ws.Range("K33").Value = 379.84
$ws.Range("L33").Value = 513.2857
$ws.Range("M33").Value = -150.84
$ws.Range("N33").Value = -971.2857
# Row 62
$ws.Range("H62").Value = 7332.6665
$ws.Range("I62").Value = 7332.6665
$ws.Range("K62").Value = 7332.6665
$ws.Range("M62").Value = -6708.6665
# Row 65
$ws.Range("H65").Value = 7332.6665
$ws.Range("I65").Value = 7332.6665
$ws.Range("K65").Value = 36663.3325
$ws.Range("M65").Value = -33543.3325
# Row 98
$ws.Range("H98").Value = 2745.2144
$ws.Range("I98").Value = 3029.6
$ws.Range("J98").Value = 1323.2858
$ws.Range("K98").Value = 3029.6
$ws.Range("L98").Value = 1323.2858
$ws.Range("M98").Value = -1531.6
$ws.Range("N98").Value = -4319.2858
# Row 100
$ws.Range("H100").Value = 2147.875
$ws.Range("I100").Value = 2172.1428
$ws.Range("K100").Value = 2172.1428
$ws.Range("M100").Value = -1631.1428
# Row 122
$ws.Range("H122").Value = 2745.2144
$ws.Range("I122").Value = 3029.6
$ws.Range("J122").Value = 1323.2858
$ws.Range("K122").Value = 9088.799999999999
$ws.Range("L122").Value = 3969.8574
$ws.Range("M122").Value = -6638.799999999999
$ws.Range("N122").Value = -8869.857400000001
# Row 132
$ws.Range("H132").Value = 6541503
$ws.Range("I132").Value = 9809149
$ws.Range("K132").Value = 29427447
$ws.Range("M132").Value = -29424917
# Row 137
$ws.Range("H137").Value = 1526.0454
$ws.Range("I137").Value = 1369.579
$ws.Range("J137").Value = 1644.96
$ws.Range("K137").Value = 4108.737
$ws.Range("L137").Value = 4934.88
$ws.Range("M137").Value = -1558.737
$ws.Range("N137").Value = -10034.88
# Row 138
$ws.Range("H138").Value = 702859.5
$ws.Range("J138").Value = 967866.5600000001
$ws.Range("L138").Value = 2903599.68
$ws.Range("N138").Value = -2913879.68
# Row 141
$ws.Range("H141").Value = 1779.1666
$ws.Range("I141").Value = 1168.75
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 3506.25
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 1673.75
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4691.385
$ws.Range("I32").Value = 3785.672
$ws.Range("J32").Value = 18503.5
$ws.Range("K32").Value = 3785.672
$ws.Range("L32").Value = 18503.5
$ws.Range("M32").Value = -3498.672
$ws.Range("N32").Value = -19077.5
# Row 61
$ws.Range("H61").Value = 914.0769
$ws.Range("I61").Value = 822.4167
$ws.Range("K61").Value = 822.4167
$ws.Range("M61").Value = -610.4167
# Row 74
$ws.Range("H74").Value = 1817.9474
$ws.Range("I74").Value = 1585.6111
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 1585.6111
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -711.6111000000001
$ws.Range("N74").Value = -7748
# Row 77
$ws.Range("H77").Value = 1817.9474
$ws.Range("I77").Value = 1585.6111
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 7928.0555
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -3560.0555
$ws.Range("N77").Value = -38736
# Row 102
$ws.Range("H102").Value = 23811980
$ws.Range("I102").Value = 33335574
$ws.Range("K102").Value = 33335574
$ws.Range("M102").Value = -33333952
# Row 119
$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676
# Row 122
$ws.Range("H122").Value = 1352.8
$ws.Range("I122").Value = 1385.1428
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 4155.428400000001
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -1705.428400000001
$ws.Range("N122").Value = -7600
# Row 124
$ws.Range("H124").Value = 17296
$ws.Range("J124").Value = 17296
$ws.Range("L124").Value = 17296
$ws.Range("N124").Value = -27116
# Row 132
$ws.Range("H132").Value = 2673.862
$ws.Range("I132").Value = 2341.8
$ws.Range("J132").Value = 4749.25
$ws.Range("K132").Value = 7025.400000000001
$ws.Range("L132").Value = 14247.75
$ws.Range("M132").Value = -4495.400000000001
$ws.Range("N132").Value = -19307.75
# Row 136
$ws.Range("H136").Value = 914.0769
$ws.Range("I136").Value = 822.4167
$ws.Range("K136").Value = 2467.2501
$ws.Range("M136").Value = 82.7498999999998

$ws = $wb.Worksheets.Item("BSM")
# Row 23
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("K23").Value = 1000
$ws.Range("M23").Value = -717
# Row 59
$ws.Range("H59").Value = 60390
$ws.Range("J59").Value = 70780
$ws.Range("L59").Value = 70780
$ws.Range("N59").Value = -72474
# Row 134
$ws.Range("H134").Value = 5187.3076
$ws.Range("I134").Value = 1036.25
$ws.Range("K134").Value = 3108.75
$ws.Range("M134").Value = -573.75

$ws = $wb.Worksheets.Item("CRP")
# Row 74
$ws.Range("H74").Value = 28333.334
$ws.Range("I74").Value = 19000
$ws.Range("J74").Value = 33000
$ws.Range("K74").Value = 19000
$ws.Range("L74").Value = 33000
$ws.Range("M74").Value = -18126
$ws.Range("N74").Value = -34748
# Row 77
$ws.Range("H77").Value = 28333.334
$ws.Range("I77").Value = 19000
$ws.Range("J77").Value = 33000
$ws.Range("K77").Value = 57000
$ws.Range("L77").Value = 99000
$ws.Range("M77").Value = -52632
$ws.Range("N77").Value = -107736
# Row 105
$ws.Range("H105").Value = 508.42856
$ws.Range("I105").Value = 459.83334
$ws.Range("K105").Value = 459.83334
$ws.Range("M105").Value = 1287.16666
# Row 122
$ws.Range("H122").Value = 1058
$ws.Range("I122").Value = 1058
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3174
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -724
# Row 132
$ws.Range("H132").Value = 1675.3684
$ws.Range("I132").Value = 1290.4814
$ws.Range("J132").Value = 2620.0908
$ws.Range("K132").Value = 3871.4442
$ws.Range("L132").Value = 7860.2724
$ws.Range("M132").Value = -1341.4442
$ws.Range("N132").Value = -12920.2724

$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 313.375
$ws.Range("I40").Value = 103
$ws.Range("J40").Value = 409
$ws.Range("K40").Value = 412
$ws.Range("L40").Value = 1636
$ws.Range("M40").Value = -343
$ws.Range("N40").Value = -1774
# Row 131
$ws.Range("H131").Value = 13336000
$ws.Range("I131").Value = 142857460
$ws.Range("J131").Value = 2908.8677
$ws.Range("K131").Value = 428572380
$ws.Range("L131").Value = 8726.6031
$ws.Range("M131").Value = -428567340
$ws.Range("N131").Value = -18806.6031

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 475.45456
$ws.Range("I102").Value = 459
$ws.Range("J102").Value = 640
$ws.Range("K102").Value = 459
$ws.Range("L102").Value = 640
$ws.Range("M102").Value = 1163
$ws.Range("N102").Value = -3884

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 35721400
$ws.Range("I122").Value = 62501500
$ws.Range("J122").Value = 14601.333
$ws.Range("K122").Value = 187504500
$ws.Range("L122").Value = 43803.999
$ws.Range("M122").Value = -187502050
$ws.Range("N122").Value = -48703.999
# Row 132
$ws.Range("H132").Value = 57315.6
$ws.Range("I132").Value = 18117.334
$ws.Range("J132").Value = 74114.86
$ws.Range("K132").Value = 54352.00199999999
$ws.Range("L132").Value = 222344.58
$ws.Range("M132").Value = -51822.00199999999
$ws.Range("N132").Value = -227404.58
# Row 136
$ws.Range("H136").Value = 1824.125
$ws.Range("I136").Value = 1398.6
$ws.Range("J136").Value = 2533.3333
$ws.Range("K136").Value = 4195.799999999999
$ws.Range("L136").Value = 7599.999899999999
$ws.Range("M136").Value = -1645.799999999999
$ws.Range("N136").Value = -12699.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 378.33334
$ws.Range("I100").Value = 378.33334
$ws.Range("K100").Value = 756.66668
$ws.Range("M100").Value = -215.66668
# Row 132
$ws.Range("H132").Value = 1309.238
$ws.Range("I132").Value = 971.17645
$ws.Range("J132").Value = 2746
$ws.Range("K132").Value = 2913.52935
$ws.Range("L132").Value = 8238
$ws.Range("M132").Value = -383.5293500000002
$ws.Range("N132").Value = -13298
